$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.175.02"
$ws.Range("E2").Value = "  -2.86%  "

$ws.Range("D3").Value = "1.849.33"
$ws.Range("E3").Value = "  -1.88%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7061"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.40%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3055"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07414"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08137"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7271"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.74%  "

$ws.Range("D13").Value = "1.834.06"
$ws.Range("E13").Value = "  -3.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.224"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "88.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.00%  "

$ws.Range("D16").Value = "29.175.21"
$ws.Range("E16").Value = "  -2.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.762"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.41%  "

$ws.Range("E19").Value = "  -3.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007651"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9994"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "2.093.92"
$ws.Range("E22").Value = "  -3.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.600"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.58%  "

$ws.Range("E25").Value = "  -2.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1456"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.973"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.399"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.529"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.44%  "

$ws.Range("E32").Value = "  -2.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.990"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05194"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.187"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.033"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7043"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.664"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.33%  "

$ws.Range("E39").Value = "  -4.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.680"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9518"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.021"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.30%  "

$ws.Range("D43").Value = "1.073.47"
$ws.Range("E43").Value = "  -1.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4293"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.744"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.76%  "

$ws.Range("D49").Value = "1.986.05"
$ws.Range("E49").Value = "  -3.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.060"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.111"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.55%  "
